# Finished state data pcp upload.
# Fill in the "api" column (C) with "yes" for the block of rows that was
# missing it (rows 67-154 on the "all" sheet), matching the rest of the
# column. The very first newly-filled cell (C67) keeps a red font, as if
# it had been typed/marked first before the rest were filled down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")

# Bulk-fill C67:C154 with "yes" (reuses the existing shared string).
$ws.Range("C67:C154").Value = "yes"

# C67 is marked in red font.
$ws.Range("C67").Font.Color = 255

# Leave the selection where the edit left off.
$ws.Range("C154").Select()
